$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (column G) values - replacing Strike# with recalculated K (strikeouts)
$kValues = @{
    2 = 2
    3 = 0
    4 = 0
    5 = 3
    6 = 0
    7 = 1
    8 = 0
    9 = 0
    10 = 2
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 2
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 0
    33 = 1
    34 = 4
    35 = 1
    36 = 0
    37 = 0
    38 = 2
    39 = 0
    40 = 2
    41 = 0
    42 = 3
    43 = 1
    44 = 0
    45 = 4
    46 = 1
    47 = 3
    48 = 5
    49 = 1
    50 = 1
    51 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

